# Weekly update: insert a new observation row for "Haba" (Terminal La
# Palmera de La Serena) as row 14, shifting all existing data rows
# (previously 14-42) down by one (to 15-43). The sheet has no table
# object, just plain cell data, so a simple row insert does the job and
# Excel keeps everything below it intact.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 14 - rows 14:42 shift down to 15:43.
$ws.Rows.Item(14).Insert()

# Populate the new row 14 with the latest weekly price observation.
$ws.Range("A14").Value = 8
$ws.Range("B14").Value = "Terminal La Palmera de La Serena"
$ws.Range("C14").Value = "Coquimbo"
$ws.Range("D14").Value = 44883
$ws.Range("E14").Value = 4
$ws.Range("F14").Value = 100112026
$ws.Range("G14").Value = "Haba"
$ws.Range("H14").Value = "Sin especificar"
$ws.Range("I14").Value = "Primera"
$ws.Range("J14").Value = 380
$ws.Range("K14").Value = 7000
$ws.Range("L14").Value = 8000
$ws.Range("M14").Value = 7500
$ws.Range("N14").Value = "$/saco 25 kilos"
$ws.Range("O14").Value = "Provincia del Elquí"
$ws.Range("P14").Value = 300
$ws.Range("Q14").Value = 25
$ws.Range("R14").Value = "Hortaliza"
